$d = $word.ActiveDocument

# --- 1. Split the "Angel:" bullet's second run so "studied for QC" becomes
#        its own run (reverting the run-merge caused by the later edit). ---
$p5 = $d.Paragraphs.Item(5)
$pText = $p5.Range.Text
$idx = $pText.IndexOf("studied for QC")
$splitPos = $p5.Range.Start + $idx + 1
$tail = $d.Range($splitPos, $p5.Range.End)
# Force Word to materialize a distinct run for the tail text by touching a
# character property and restoring it to its original value.
$tail.Font.Color = 1
$tail.Font.Color = 0

# --- 2. Remove everything between the "Login" paragraph and the trailing
#        empty paragraph (the standup notes that were merged back in). ---
$loginPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -eq "Login" + [char]13) {
        $loginPara = $cand
    }
}

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$delRange = $d.Range($loginPara.Range.End, $lastPara.Range.Start)
$delRange.Delete()

# --- 3. Drop the now-unused "List Paragraph" style definition. ---
$listParaStyle = $d.Styles("ListParagraph")
$listParaStyle.Delete()
